$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the very end of the paragraph
# (it marks the author's last edit position). After the edit it should sit
# right between the new "con " and "CI. ..." text, so drop it now and
# re-create it later once that text exists.
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
}

# Replace the old clause ("más controles por consulta externa desde las ")
# with the new prose, keeping the trailing "desde las " (which feeds the
# start_time MERGEFIELD right after it) intact.
$r = $d.Content
$null = $r.Find.Execute(
    "más controles por consulta externa desde las ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "acude en compañía de su madre, la ... con CI. ... desde las ",
    2)

# Turn the single leading space before that clause into ", ".
$r2 = $d.Content
$null = $r2.Find.Execute(
    "medicamentoso ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "medicamentoso, ",
    2)

# Re-insert the "_GoBack" bookmark right after "con " (i.e. right before
# "CI. ...").
$r3 = $d.Content
$found3 = $r3.Find.Execute(
    "la ... con ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "")
if ($found3) {
    $r3.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $r3)
}
